# Update: Read only today's emails and restructure template layout
# The log previously accumulated two email entries (rows 2 and 3). The
# template now only keeps a single, freshly downloaded entry (row 2),
# reflecting today's email results, and the old second entry is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 with the latest email/download record ---
$ws.Range("A2").Value2 = 'Project Documents: TSE Enquiry and Quotation for Review'
$ws.Range("B2").Value2 = '197c45252cdeeb84'
$ws.Range("C2").Value2 = '197c45252cdeeb84'
$ws.Range("D2").Value2 = 'Kumar <kum4r18@gmail.com>'

$ws.Range("E2").Value2 = 45839.54337962963
$ws.Range("F2").Value2 = 45839.54337962963
$ws.Range("G2").Value2 = 45839.54337962963
$ws.Range("H2").Value2 = 45839.54337962963

$ws.Range("I2").Value2 = 2
$ws.Range("J2").Value2 = 0

$ws.Range("K2").Value2 = 'quotation_20250701_130227_TSE Enquiry form-Mega 80S.pdf, quotation_20250701_130228_QU-IMM-Vi-42025-00169-1-28-04-2025-STEER ENGG(REV).pdf'
$ws.Range("L2").Value2 = '/mnt/c/Users/Imran/OneDrive - Ahana Systems and Solutions (P) Ltd/Desktop/Demo/steer_document_processing_poc/demo_app/backend/Agent_AI/download_email/quotation_20250701_130227_TSE Enquiry form-Mega 80S.pdf, /mnt/c/Users/Imran/OneDrive - Ahana Systems and Solutions (P) Ltd/Desktop/Demo/steer_document_processing_poc/demo_app/backend/Agent_AI/download_email/quotation_20250701_130228_QU-IMM-Vi-42025-00169-1-28-04-2025-STEER ENGG(REV).pdf'
$ws.Range("N2").Value2 = '/mnt/c/Users/Imran/OneDrive - Ahana Systems and Solutions (P) Ltd/Desktop/Demo/steer_document_processing_poc/demo_app/backend/Agent_AI/result_json/quotation_20250701_130227_TSE Enquiry form-Mega 80S.json,/mnt/c/Users/Imran/OneDrive - Ahana Systems and Solutions (P) Ltd/Desktop/Demo/steer_document_processing_poc/demo_app/backend/Agent_AI/result_json/quotation_20250701_130228_QU-IMM-Vi-42025-00169-1-28-04-2025-STEER ENGG(REV).json'

$ws.Range("O2").Value2 = 'pending,pending'
$ws.Range("P2").Value2 = 'c2aaf7d4be1f9bd055873ae4fea6e206'
$ws.Range("Q2").Value2 = '179121bfa9569741a3714028a4e7b3f25d0705f5d86e409039971ff6a0aeabf2, 206fabc798763a2e17c5de22364ed1088884d2926bcf51153b3af314babd7585'
$ws.Range("R2").Value2 = 'TSE Enquiry form-Mega 80S.pdf_197c45252cdeeb84_179121bfa9569741, QU-IMM-Vi-42025-00169-1-28-04-2025-STEER ENGG(REV).pdf_197c45252cdeeb84_206fabc798763a2e'

$ws.Range("T2").Value2 = 'QUOTATION'
$ws.Range("V2").Value2 = 'completed,completed'

# --- Remove the now-obsolete second email entry (old row 3) ---
$ws.Rows.Item(3).Delete()
